# Edit script: applies textual + structural changes described by the diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "7 distinct measures" -> "8 distinct measures"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The Border Crossing/Entry dataset contains 364,510 records, pertaining to 7 distinct measures",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The Border Crossing/Entry dataset contains 364,510 records, pertaining to 8 distinct measures",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Flesh out the "4.Discussion" summary line
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "4.Discussion: periodicity of the port cities; pedestrian v port traffic;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "4.Discussion: The border crossing dataset presents various interesting facets, to include periodicity of the port cities, differences between northern and southern borders, anomolies, and future forecasts.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Drop the "First: " lead-in on the periodicity bullet
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "First: Every border crossing city experiences periodicity",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Every border crossing city experiences periodicity",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Drop the "Second: " lead-in and extend the footprint bullet
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Second: There exists a significantly larger traffic footprint at the U.S.-Mexico border, than the U.S.-Canada border. Here is where I talk Andra" + [char]8217 + "s pivot data" + [char]8230,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "There exists a significantly larger traffic footprint at the U.S.-Mexico border, than the U.S.-Canada border. Although inference is not made with regards to reasons, in the largest metric, Personal Passenger Vehicles, the largest yearly value for the US-Mexico border is nearly doube that of the US-Canada border.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Drop the "Third: " lead-in
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Third: Anomolies exists within the dataset",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Anomolies exists within the dataset",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Drop the "Fourth: " lead-in
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Fourth: Forecasts for any given port",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Forecasts for any given port",
    2) | Out-Null

Write-Output "Text replacements done"
